$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("WMT_Extract")
$ws1.Copy($null, $ws1)
$newSheet = $wb.Worksheets.Item("WMT_Extract (2)")
$newSheet.Name = "WMT_Extract_Filtered"
$ws1.Rows("1:3").Select()
